$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-10) so the shared-string table rebuilds cleanly
$ws.Range("A2:T10").ClearContents()

# Set string columns in column-major order (A, then B, then C, then D across rows 2-13)
# so the rebuilt shared-string table orders new entries as: ECs, Inflammatory-Mac, Neutrophils,
# Resolving-Mac, Ccl12, Ackr4, FAPs, MuSCs (matching the target workbook).
$ws.Range("A2").Value2 = "ECs"
$ws.Range("A3").Value2 = "ECs"
$ws.Range("A4").Value2 = "ECs"
$ws.Range("A5").Value2 = "Inflammatory-Mac"
$ws.Range("A6").Value2 = "Inflammatory-Mac"
$ws.Range("A7").Value2 = "Inflammatory-Mac"
$ws.Range("A8").Value2 = "Neutrophils"
$ws.Range("A9").Value2 = "Neutrophils"
$ws.Range("A10").Value2 = "Neutrophils"
$ws.Range("A11").Value2 = "Resolving-Mac"
$ws.Range("A12").Value2 = "Resolving-Mac"
$ws.Range("A13").Value2 = "Resolving-Mac"

$ws.Range("B2").Value2 = "Ccl12"
$ws.Range("B3").Value2 = "Ccl12"
$ws.Range("B4").Value2 = "Ccl12"
$ws.Range("B5").Value2 = "Ccl12"
$ws.Range("B6").Value2 = "Ccl12"
$ws.Range("B7").Value2 = "Ccl12"
$ws.Range("B8").Value2 = "Ccl12"
$ws.Range("B9").Value2 = "Ccl12"
$ws.Range("B10").Value2 = "Ccl12"
$ws.Range("B11").Value2 = "Ccl12"
$ws.Range("B12").Value2 = "Ccl12"
$ws.Range("B13").Value2 = "Ccl12"

$ws.Range("C2").Value2 = "Ackr4"
$ws.Range("C3").Value2 = "Ackr4"
$ws.Range("C4").Value2 = "Ackr4"
$ws.Range("C5").Value2 = "Ackr4"
$ws.Range("C6").Value2 = "Ackr4"
$ws.Range("C7").Value2 = "Ackr4"
$ws.Range("C8").Value2 = "Ackr4"
$ws.Range("C9").Value2 = "Ackr4"
$ws.Range("C10").Value2 = "Ackr4"
$ws.Range("C11").Value2 = "Ackr4"
$ws.Range("C12").Value2 = "Ackr4"
$ws.Range("C13").Value2 = "Ackr4"

$ws.Range("D2").Value2 = "ECs"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("D10").Value2 = "MuSCs"
$ws.Range("D11").Value2 = "ECs"
$ws.Range("D12").Value2 = "FAPs"
$ws.Range("D13").Value2 = "MuSCs"

# Set numeric columns
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 1.005755333333333
$ws.Range("H2").Value2 = 3.017266
$ws.Range("I2").Value2 = 0.01048729000197281
$ws.Range("J2").Value2 = 0.01048729000197281
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.180428
$ws.Range("N2").Value2 = 0.541284
$ws.Range("O2").Value2 = 0.6724860231084607
$ws.Range("P2").Value2 = 0.6724860231084607
$ws.Range("Q2").Value2 = 0.1814664232826667
$ws.Range("R2").Value2 = 1.633197809544
$ws.Range("S2").Value2 = 0.007052555946611814
$ws.Range("T2").Value2 = 0.007052555946611814

$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 1.005755333333333
$ws.Range("H3").Value2 = 3.017266
$ws.Range("I3").Value2 = 0.01048729000197281
$ws.Range("J3").Value2 = 0.01048729000197281
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.01727566666666666
$ws.Range("N3").Value2 = 0.051827
$ws.Range("O3").Value2 = 0.06438936513852653
$ws.Range("P3").Value2 = 0.06438936513852653
$ws.Range("Q3").Value2 = 0.01737509388688889
$ws.Range("R3").Value2 = 0.156375844982
$ws.Range("S3").Value2 = 0.0006752699452506456
$ws.Range("T3").Value2 = 0.0006752699452506456

$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 1.005755333333333
$ws.Range("H4").Value2 = 3.017266
$ws.Range("I4").Value2 = 0.01048729000197281
$ws.Range("J4").Value2 = 0.01048729000197281
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 0.07059633333333333
$ws.Range("N4").Value2 = 0.211789
$ws.Range("O4").Value2 = 0.2631246117530128
$ws.Range("P4").Value2 = 0.2631246117530128
$ws.Range("Q4").Value2 = 0.07100263876377777
$ws.Range("R4").Value2 = 0.6390237488739999
$ws.Range("S4").Value2 = 0.002759464110110348
$ws.Range("T4").Value2 = 0.002759464110110348

$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 38.755371
$ws.Range("H5").Value2 = 116.266113
$ws.Range("I5").Value2 = 0.4041130097356814
$ws.Range("J5").Value2 = 0.4041130097356814
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.180428
$ws.Range("N5").Value2 = 0.541284
$ws.Range("O5").Value2 = 0.6724860231084607
$ws.Range("P5").Value2 = 0.6724860231084607
$ws.Range("Q5").Value2 = 6.992554078787999
$ws.Range("R5").Value2 = 62.93298670909199
$ws.Range("S5").Value2 = 0.2717603508035391
$ws.Range("T5").Value2 = 0.2717603508035391

$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 38.755371
$ws.Range("H6").Value2 = 116.266113
$ws.Range("I6").Value2 = 0.4041130097356814
$ws.Range("J6").Value2 = 0.4041130097356814
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.01727566666666666
$ws.Range("N6").Value2 = 0.051827
$ws.Range("O6").Value2 = 0.06438936513852653
$ws.Range("P6").Value2 = 0.06438936513852653
$ws.Range("Q6").Value2 = 0.6695248709389999
$ws.Range("R6").Value2 = 6.025723838450999
$ws.Range("S6").Value2 = 0.02602058014109971
$ws.Range("T6").Value2 = 0.02602058014109971

$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 38.755371
$ws.Range("H7").Value2 = 116.266113
$ws.Range("I7").Value2 = 0.4041130097356814
$ws.Range("J7").Value2 = 0.4041130097356814
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 0.6666666666666666
$ws.Range("M7").Value2 = 0.07059633333333333
$ws.Range("N7").Value2 = 0.211789
$ws.Range("O7").Value2 = 0.2631246117530128
$ws.Range("P7").Value2 = 0.2631246117530128
$ws.Range("Q7").Value2 = 2.735987089573
$ws.Range("R7").Value2 = 24.623883806157
$ws.Range("S7").Value2 = 0.1063320787910427
$ws.Range("T7").Value2 = 0.1063320787910427

$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 7.988471999999999
$ws.Range("H8").Value2 = 23.965416
$ws.Range("I8").Value2 = 0.08329801469605898
$ws.Range("J8").Value2 = 0.08329801469605898
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 0.180428
$ws.Range("N8").Value2 = 0.541284
$ws.Range("O8").Value2 = 0.6724860231084607
$ws.Range("P8").Value2 = 0.6724860231084607
$ws.Range("Q8").Value2 = 1.441344026016
$ws.Range("R8").Value2 = 12.972096234144
$ws.Range("S8").Value2 = 0.05601675063578282
$ws.Range("T8").Value2 = 0.05601675063578282

$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 7.988471999999999
$ws.Range("H9").Value2 = 23.965416
$ws.Range("I9").Value2 = 0.08329801469605898
$ws.Range("J9").Value2 = 0.08329801469605898
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.01727566666666666
$ws.Range("N9").Value2 = 0.051827
$ws.Range("O9").Value2 = 0.06438936513852653
$ws.Range("P9").Value2 = 0.06438936513852653
$ws.Range("Q9").Value2 = 0.138006179448
$ws.Range("R9").Value2 = 1.242055615032
$ws.Range("S9").Value2 = 0.00536350628357889
$ws.Range("T9").Value2 = 0.00536350628357889

$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 7.988471999999999
$ws.Range("H10").Value2 = 23.965416
$ws.Range("I10").Value2 = 0.08329801469605898
$ws.Range("J10").Value2 = 0.08329801469605898
$ws.Range("K10").Value2 = 2
$ws.Range("L10").Value2 = 0.6666666666666666
$ws.Range("M10").Value2 = 0.07059633333333333
$ws.Range("N10").Value2 = 0.211789
$ws.Range("O10").Value2 = 0.2631246117530128
$ws.Range("P10").Value2 = 0.2631246117530128
$ws.Range("Q10").Value2 = 0.5639568321359999
$ws.Range("R10").Value2 = 5.075611489223999
$ws.Range("S10").Value2 = 0.02191775777669728
$ws.Range("T10").Value2 = 0.02191775777669728

$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 48.15271133333334
$ws.Range("H11").Value2 = 144.458134
$ws.Range("I11").Value2 = 0.5021016855662869
$ws.Range("J11").Value2 = 0.5021016855662868
$ws.Range("K11").Value2 = 2
$ws.Range("L11").Value2 = 0.6666666666666666
$ws.Range("M11").Value2 = 0.180428
$ws.Range("N11").Value2 = 0.541284
$ws.Range("O11").Value2 = 0.6724860231084607
$ws.Range("P11").Value2 = 0.6724860231084607
$ws.Range("Q11").Value2 = 8.688097400450667
$ws.Range("R11").Value2 = 78.192876604056
$ws.Range("S11").Value2 = 0.3376563657225271
$ws.Range("T11").Value2 = 0.337656365722527

$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 48.15271133333334
$ws.Range("H12").Value2 = 144.458134
$ws.Range("I12").Value2 = 0.5021016855662869
$ws.Range("J12").Value2 = 0.5021016855662868
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Value2 = 0.3333333333333333
$ws.Range("M12").Value2 = 0.01727566666666666
$ws.Range("N12").Value2 = 0.051827
$ws.Range("O12").Value2 = 0.06438936513852653
$ws.Range("P12").Value2 = 0.06438936513852653
$ws.Range("Q12").Value2 = 0.8318701900908888
$ws.Range("R12").Value2 = 7.486831710818
$ws.Range("S12").Value2 = 0.03233000876859728
$ws.Range("T12").Value2 = 0.03233000876859727

$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 48.15271133333334
$ws.Range("H13").Value2 = 144.458134
$ws.Range("I13").Value2 = 0.5021016855662869
$ws.Range("J13").Value2 = 0.5021016855662868
$ws.Range("K13").Value2 = 2
$ws.Range("L13").Value2 = 0.6666666666666666
$ws.Range("M13").Value2 = 0.07059633333333333
$ws.Range("N13").Value2 = 0.211789
$ws.Range("O13").Value2 = 0.2631246117530128
$ws.Range("P13").Value2 = 0.2631246117530128
$ws.Range("Q13").Value2 = 3.399404860191778
$ws.Range("R13").Value2 = 30.594643741726
$ws.Range("S13").Value2 = 0.1321153110751626
$ws.Range("T13").Value2 = 0.1321153110751625

